$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108, shifting the existing row 108 (and all
# rows below it) down by one. This preserves all existing data/formatting
# and creates space for the new weekly record.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new weekly price record.
$ws.Range("A108").Value = 8
$ws.Range("B108").Value = "Terminal La Palmera de La Serena"
$ws.Range("C108").Value = "Coquimbo"
$ws.Range("D108").Value = 44960
$ws.Range("E108").Value = 4
$ws.Range("F108").Value = 100112052
$ws.Range("G108").Value = "Albahaca"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 960
$ws.Range("K108").Value = 5000
$ws.Range("L108").Value = 6000
$ws.Range("M108").Value = 5500
$ws.Range("N108").Value = "`$/docena de matas"
$ws.Range("O108").Value = "Provincia del Elquí"
$ws.Range("P108").Value = 917
$ws.Range("Q108").Value = 6
$ws.Range("R108").Value = "Hortaliza"
